$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updated values (columns B..Q; A, O and R are unchanged)
$ws.Range("B2").Value = 11554617.8962204
$ws.Range("C2").Value = -268947467.421113
$ws.Range("D2").Value = 150420.128107754
$ws.Range("E2").Value = 113487.224209943
$ws.Range("F2").Value = 278274299.226709
$ws.Range("G2").Value = 278387786.450919
$ws.Range("H2").Value = 6329.2755529044
$ws.Range("I2").Value = 0.0407523589586967
$ws.Range("J2").Value = 65.8656047451074
$ws.Range("K2").Value = 5993.77003180478
$ws.Range("L2").Value = 18065.8656047451
$ws.Range("M2").Value = 27093.9020234045
$ws.Range("N2").Value = 280171313.839889
$ws.Range("P2").Value = 1361104165277.77
$ws.Range("Q2").Value = 1361384336591.61

# Row 3 updated values (columns B..Q; A, O and R are unchanged)
$ws.Range("B3").Value = 4037218.5925632
$ws.Range("C3").Value = -153492987.85454
$ws.Range("D3").Value = 103281.492049934
$ws.Range("E3").Value = 78385.4398651211
$ws.Range("F3").Value = 156770400.673036
$ws.Range("G3").Value = 156848786.112901
$ws.Range("H3").Value = 6329.2755529044
$ws.Range("I3").Value = 0.0407523589586967
$ws.Range("J3").Value = 65.8656047451074
$ws.Range("K3").Value = 5993.77003180478
$ws.Range("L3").Value = 18065.8656047451
$ws.Range("M3").Value = 27093.9020234045
$ws.Range("N3").Value = 93485734.7981108
$ws.Range("P3").Value = 766800415826.048
$ws.Range("Q3").Value = 766893901560.847
